$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new bug tracking row (row 3), reusing number formats from row 2
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A3").Value = 42263.895833333336
$ws.Range("B3").Value = "[Html] Footer won't stay down"
$ws.Range("C3").Value = "Henry"
$ws.Range("D3").Value = 42263
$ws.Range("E3").Value = "Web App"

# Update selection to match final state
$ws.Range("B4").Select()
